$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 100
$ws1.Range("F9").Value = 8637
$ws1.Range("F19").Value = 62
$ws1.Range("F21").Value = 997

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 100
$ws4.Range("F11").Value = 8637
$ws4.Range("F21").Value = 62
$ws4.Range("F23").Value = 997
